$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 384
$ws1.Range("G2").Value = "不可售"
$ws1.Range("F3").Value = 844
$ws1.Range("G3").Value = "不可售"
$ws1.Range("F5").Value = 1047
$ws1.Range("F6").Value = 2437
$ws1.Range("F7").Value = 206

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 384
$ws4.Range("G2").Value = "不可售"
$ws4.Range("F3").Value = 844
$ws4.Range("G3").Value = "不可售"
$ws4.Range("F7").Value = 1047
$ws4.Range("F8").Value = 2437
$ws4.Range("F10").Value = 206
